# Remove the trailing "Ver no Jupiter ..." / copyright footer block that
# followed the last LOT course line, along with the blank paragraph that
# separated it from that line, leaving the blank paragraph that precedes
# the final page-break paragraph untouched.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$jupiterIdx = -1
$contactIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Ver no Jupiter") { $jupiterIdx = $i }
    if ($t -match "Contact:\s*luizeleno") { $contactIdx = $i }
}

if ($jupiterIdx -gt 0 -and $contactIdx -ge $jupiterIdx) {
    # The blank paragraph right before the "Ver no Jupiter" paragraph is
    # also removed (it only separated the LOT line from the footer block).
    $startIdx = $jupiterIdx - 1
    $startPara = $d.Paragraphs.Item($startIdx)
    $endPara = $d.Paragraphs.Item($contactIdx)

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
